$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels to include units (commented code and added units)
$ws.Range("B1").Value = "amplitude (counts)"
$ws.Range("C1").Value = "center (degrees)"
$ws.Range("D1").Value = "sigma (width in degrees)"
$ws.Range("E1").Value = "FWHM (degrees)"
$ws.Range("G1").Value = "d (Angstroms)"
